# "chachi update for new year"
# Target sheet: "Misc" (sheet4.xml) — already the ActiveSheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 blank rows at row 697, pushing the existing rows 697-748 down
#     to 701-752 (this recreates the row-shift seen from row 698 onward in
#     the diff; new data then lands in the freshly-opened 697/698 rows). ---
$ws.Rows("697:700").Insert()

# New rows 697-698 (two new card entries, player name re-uses existing
# shared strings; image URLs are brand-new shared strings — order below
# matches the shared-string append order required by the diff).
$ws.Range("B697").Value = "https://4.bp.blogspot.com/-jeTScXLhNBE/WGkyY0JZ3mI/AAAAAAAAe24/dPKNVpHsd5opYzEIg1-_mqNOGRXiwUdtQCLcB/s1600/1967%2BTopps%2B%2523560%2BBunning.jpg"
$ws.Range("A697").Value = "Jim Bunning"

$ws.Range("B698").Value = "https://1.bp.blogspot.com/-meu1Keyba7k/WGkyYxHHlvI/AAAAAAAAe28/by7E7RVUoKk15FsArlMH108aIv9Q_A1jwCLcB/s1600/1967%2BTopps%2B%2523595%2BRojas.jpg"
$ws.Range("A698").Value = "Cookie Rojas"

# Fill in existing gap rows with new card entries (no row shifting needed,
# these row numbers were already unused/empty in the sheet).
$ws.Range("A517").Value = "Ozzie Virgil"
$ws.Range("B517").Value = "https://1.bp.blogspot.com/-jxjf7CvIeQA/WGHqfuRUWRI/AAAAAAABq3M/jEXUzUPgsyQ29OSZncHZs3dHy0nzQQQrgCLcB/s1600/Virgil%2B85T.jpg"

$ws.Range("B484").Value = "https://1.bp.blogspot.com/-6bRvw9eeEyY/WF7mdXfRLGI/AAAAAAAAgtY/RKDhzd3NME85FMII5tLLzjtLroC31hThgCLcB/s1600/pippen9697hoops.jpg"
$ws.Range("A484").Value = "Scottie Pippen  1996-97 "

$ws.Range("A326").Value = "Pete Rose 1985 Renata Galasso"
$ws.Range("B326").Value = "https://2.bp.blogspot.com/-JB8nIqBe16s/WGVV_aF8Q7I/AAAAAAAAcXU/5GEgYYG5uVgor_lK_jB2NN5a5r27_YBlwCLcB/s1600/Rose%2BRenata%2BGalasso%2BPete%2BRose%2BF.jpg"

$ws.Range("A473").Value = "Steve Carlton Giants"
$ws.Range("B473").Value = "https://1.bp.blogspot.com/-3Bi43Defgv8/WGQlGlgxFsI/AAAAAAAAQ8c/KpoLGOJ6TIwNqZ9n9s8tcycA2G1g9hLqwCLcB/s1600/FREE%2BCARDS%2B009%2B3.jpg"

# Match the final selection state recorded in the workbook view.
$ws.Range("B473").Select()
